$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "Matteo Alberti"
$ws.Range("B44").Value = "Thomas Debiasi | Mai una gioia"
$ws.Range("C44").Value = "Nadir  chtioui | Mai una gioia"
$ws.Range("D44").Value = "Michael Bertè  | A.C.DENTI"
$ws.Range("E44").Value = "Andreas Galli | SdrumALA"
$ws.Range("F44").Value = "Lorenzo Zuani | I Magnifici"
